$wb = $excel.ActiveWorkbook

$reviewWs  = $wb.Worksheets.Item("REVIEW-SHEET")
$versionWs = $wb.Worksheets.Item("VERSION-HISTORY")

# --- VERSION-HISTORY: append the v1.7 entry ---------------------------------
# Copy the formatting of the last existing row (row 8) down into the new row
# (row 9) so the new row matches the look of the table (fills, borders,
# alignment, number format, etc.), then fill in the new values.
$srcRow = $versionWs.Range("A8:D8")
$dstRow = $versionWs.Range("A9:D9")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$versionWs.Cells.Item(9, 1).Value = "v1.7"
$versionWs.Cells.Item(9, 2).Value = "Mahmoud Abdelmageed"
$versionWs.Cells.Item(9, 3).Value = "Changed reviewer verification for Registration and System Constrains"
$versionWs.Cells.Item(9, 4).Value = 45766

# Match the wrapped-text row height used for the other long entries.
$versionWs.Rows.Item(9).RowHeight = 37.2

# --- Scroll REVIEW-SHEET's view back to the top (F1) ------------------------
$reviewWs.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1

# --- Make VERSION-HISTORY the active/selected sheet & cell -----------------
$versionWs.Activate()
$versionWs.Range("C9").Select()

Write-Host "Applied v1.7 version-history update"
